$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.958.07'
$ws.Range('E2').Value = '  -3.28%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.859.93'
$ws.Range('E3').Value = '  -2.48%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '317.74'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.000'
$ws.Range('E6').Value = '  -0.12%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4358'
$ws.Range('E7').Value = '  -4.85%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3696'
$ws.Range('E8').Value = '  -3.05%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07489'
$ws.Range('E9').Value = '  -3.05%  '
$ws.Range('E10').Value = '  -3.98%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '21.29'
$ws.Range('E11').Value = '  -3.51%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.878.45'
$ws.Range('E12').Value = '  -2.66%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.720'
$ws.Range('E13').Value = '  -3.16%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.435'
$ws.Range('E14').Value = '  -4.26%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.06847'
$ws.Range('E15').Value = '  -3.38%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.002'
$ws.Range('E16').Value = '  -0.12%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '81.60'
$ws.Range('E17').Value = '  -2.54%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000009023'
$ws.Range('E18').Value = '  -4.54%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.9998'
$ws.Range('E19').Value = '  -0.17%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '15.94'
$ws.Range('E20').Value = '  -4.00%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '27.930.21'
$ws.Range('E21').Value = '  -3.39%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.111'
$ws.Range('E22').Value = '  -3.86%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '11.06'
$ws.Range('E23').Value = '  +1.22%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.083.04'
$ws.Range('E24').Value = '  -3.40%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.006'
$ws.Range('E25').Value = '  -4.34%  '
$ws.Range('E26').Value = '  -2.62%  '
$ws.Range('E27').Value = '  -3.27%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '5.383'
$ws.Range('E28').Value = '  -4.97%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '113.63'
$ws.Range('E29').Value = '  -3.30%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.732'
$ws.Range('E30').Value = '  -7.44%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08982'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.8085'
$ws.Range('E32').Value = '  -6.15%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.821'
$ws.Range('E33').Value = '  -5.20%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.174'
$ws.Range('E34').Value = '  -5.89%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.945'
$ws.Range('E35').Value = '  -3.93%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.000'
$ws.Range('E36').Value = '  -0.12%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.05476'
$ws.Range('E37').Value = '  -4.03%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.112'
$ws.Range('E38').Value = '  -3.99%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01975'
$ws.Range('E39').Value = '  -3.20%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.924'
$ws.Range('E40').Value = '  +1.27%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.5251'
$ws.Range('E41').Value = '  -4.30%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '7.004'
$ws.Range('E42').Value = '  -5.59%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.1684'
$ws.Range('E43').Value = '  -3.74%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.782'
$ws.Range('E44').Value = '  -5.75%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.06794'
$ws.Range('E45').Value = '  -1.26%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4887'
$ws.Range('E46').Value = '  -5.33%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.66'
$ws.Range('E47').Value = '  -4.90%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '106.19'
$ws.Range('B49').Value = 'NEARProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.679'
$ws.Range('E49').Value = '  -5.41%  '
$ws.Range('B50').Value = 'RenderToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.911'
$ws.Range('E50').Value = '  -10.26%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.9996'
$ws.Range('E51').Value = '  -0.16%  '
